$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 3-10 and 12-13 (column A names with their B:F metric values) are being
# reordered/sorted; row 11 ("cem-vae") is unchanged. Each algorithm keeps its own
# original data values - only its row position changes.

$data = @{
    "cem"          = @(1,                    5.77, 0.9154653998284462, 0.8401729275609791, 0.8349427956342698)
    "wachter"      = @(0.5600000000000001,   3.75, 0.4935082787960795, 0.3940633180178528, 0.3940707594163054)
    "face-epsilon" = @(0.95,                 7.47, 3.952951915305309,  3.475175654078315,  0.9432876712328766)
    "cchvae"       = @(1,                    9.59, 4.288725759924322,  3.718921594863702,  0.9920547244320177)
    "ar"           = @(0.26,                 0.55, 0.4004687407569555, 0.4000020792020286, 0.26)
    "cruds"        = @(0.45,                 4.17, 2.267905295237643,  1.911294723527047,  0.4478978382746379)
    "face-knn"     = @(1,                    8.08, 4.312361890239652,  3.830199795020054,  1)
    "clue"         = @(1,                    8.65, 3.454601607651509,  2.811658536566437,  0.9927476534247398)
    "gs"           = @(1,                    6.17, 1.366335199349296,  1.184234398993123,  0.8371560574946735)
    "dice"         = @(1,                    2.39, 1.748663095629701,  1.521672553600215,  0.923548830246621)
}

$rowOrder = @{
    3  = "cem"
    4  = "wachter"
    5  = "face-epsilon"
    6  = "cchvae"
    7  = "ar"
    8  = "cruds"
    9  = "face-knn"
    10 = "clue"
    12 = "gs"
    13 = "dice"
}

foreach ($r in $rowOrder.Keys) {
    $name = $rowOrder[$r]
    $vals = $data[$name]

    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
}
